$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 154, shifting existing rows 154..216 down to 155..217
$ws.Rows.Item(154).Insert()

# Populate the newly inserted row 154 with the new data record
$ws.Range("A154").Value = 3
$ws.Range("B154").Value = "Femacal de La Calera"
$ws.Range("C154").Value = "Coquimbo"
$ws.Range("D154").Value = 44784
$ws.Range("D154").NumberFormat = $ws.Range("D155").NumberFormat
$ws.Range("E154").Value = 5
$ws.Range("F154").Value = 100112010
$ws.Range("G154").Value = "Achicoria"
$ws.Range("H154").Value = "Sin especificar"
$ws.Range("I154").Value = "Primera"
$ws.Range("J154").Value = 50
$ws.Range("K154").Value = 7000
$ws.Range("L154").Value = 7000
$ws.Range("M154").Value = 7000
$ws.Range("N154").Value = "$/caja 16 unidades"
$ws.Range("O154").Value = "Provincia de Quillota"
$ws.Range("P154").Value = 438
$ws.Range("Q154").Value = 16
$ws.Range("R154").Value = "Hortaliza"
